{"js": "const replacements = [\n  [\"720\u00d78=5760\", \"662\u00d73=1986\"],\n  [\"774\u00d73=2322\", \"144\u00d78=1152\"],\n  [\"806\u00d72=1612\", \"775\u00d72=1550\"],\n  [\"258\u00d79=2322\", \"445\u00d74=1780\"],\n  [\"978\u00d78=7824\", \"311\u00d78=2488\"],\n  [\"172\u00d74=688\", \"149\u00d77=1043\"],\n  [\"579\u00d78=4632\", \"580\u00d72=1160\"],\n  [\"198\u00d76=1188\", \"800\u00d75=4000\"],\n  [\"218\u00d78=1744\", \"295\u00d72=590\"],\n  [\"246\u00d74=984\", \"966\u00d76=5796\"],\n  [\"166\u00d73=498\", \"531\u00d77=3717\"],\n  [\"169\u00d74=676\", \"951\u00d78=7608\"],\n  [\"438\u00d75=2190\", \"787\u00d73=2361\"],\n  [\"470\u00d76=2820\", \"507\u00d72=1014\"],\n  [\"271\u00d73=813\", \"776\u00d72=1552\"],\n  [\"442\u00d79=3978\", \"899\u00d75=4495\"],\n  [\"507\u00d77=3549\", \"999\u00d73=2997\"],\n  [\"126\u00d72=252\", \"119\u00d79=1071\"],\n  [\"861\u00d73=2583\", \"419\u00d72=838\"],\n  [\"976\u00d78=7808\", \"113\u00d79=1017\"],\n  [\"649\u00d74=2596\", \"687\u00d79=6183\"],\n  [\"187\u00d76=1122\", \"509\u00d75=2545\"],\n  [\"853\u00d75=4265\", \"698\u00d78=5584\"],\n  [\"475\u00d79=4275\", \"757\u00d72=1514\"],\n  [\"717\u00d72=1434\", \"267\u00d73=801\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"720\u00d78=5760\"; New = \"662\u00d73=1986\" },\n    @{ Old = \"774\u00d73=2322\"; New = \"144\u00d78=1152\" },\n    @{ Old = \"806\u00d72=1612\"; New = \"775\u00d72=1550\" },\n    @{ Old = \"258\u00d79=2322\"; New = \"445\u00d74=1780\" },\n    @{ Old = \"978\u00d78=7824\"; New = \"311\u00d78=2488\" },\n    @{ Old = \"172\u00d74=688\";  New = \"149\u00d77=1043\" },\n    @{ Old = \"579\u00d78=4632\"; New = \"580\u00d72=1160\" },\n    @{ Old = \"198\u00d76=1188\"; New = \"800\u00d75=4000\" },\n    @{ Old = \"218\u00d78=1744\"; New = \"295\u00d72=590\" },\n    @{ Old = \"246\u00d74=984\";  New = \"966\u00d76=5796\" },\n    @{ Old = \"166\u00d73=498\";  New = \"531\u00d77=3717\" },\n    @{ Old = \"169\u00d74=676\";  New = \"951\u00d78=7608\" },\n    @{ Old = \"438\u00d75=2190\"; New = \"787\u00d73=2361\" },\n    @{ Old = \"470\u00d76=2820\"; New = \"507\u00d72=1014\" },\n    @{ Old = \"271\u00d73=813\";  New = \"776\u00d72=1552\" },\n    @{ Old = \"442\u00d79=3978\"; New = \"899\u00d75=4495\" },\n    @{ Old = \"507\u00d77=3549\"; New = \"999\u00d73=2997\" },\n    @{ Old = \"126\u00d72=252\";  New = \"119\u00d79=1071\" },\n    @{ Old = \"861\u00d73=2583\"; New = \"419\u00d72=838\" },\n    @{ Old = \"976\u00d78=7808\"; New = \"113\u00d79=1017\" },\n    @{ Old = \"649\u00d74=2596\"; New = \"687\u00d79=6183\" },\n    @{ Old = \"187\u00d76=1122\"; New = \"509\u00d75=2545\" },\n    @{ Old = \"853\u00d75=4265\"; New = \"698\u00d78=5584\" },\n    @{ Old = \"475\u00d79=4275\"; New = \"757\u00d72=1514\" },\n    @{ Old = \"717\u00d72=1434\"; New = \"267\u00d73=801\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute(\n        $r.Old,    # FindText\n        $true,     # MatchCase\n        $false,    # MatchWholeWord\n        $false,    # MatchWildcards\n        $false,    # MatchSoundsLike\n        $false,    # MatchAllWordForms\n        $true,     # Forward\n        1,         # Wrap = wdFindContinue\n        $false,    # Format\n        $r.New,    # ReplaceWith\n        2          # Replace = wdReplaceAll\n    ) | Out-Null\n}\n"}
